$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 663.55554
$ws.Range("I12").Value = 599.5714
$ws.Range("J12").Value = 887.5
$ws.Range("K12").Value = 599.5714
$ws.Range("L12").Value = 887.5
$ws.Range("M12").Value = -429.5714
$ws.Range("N12").Value = -1227.5
$ws.Range("H32").Value = 799.9375
$ws.Range("J32").Value = 988.7778
$ws.Range("L32").Value = 988.7778
$ws.Range("N32").Value = -1640.7778
$ws.Range("H41").Value = 1026.8334
$ws.Range("I41").Value = 640.25
$ws.Range("K41").Value = 640.25
$ws.Range("M41").Value = -200.25
$ws.Range("H42").Value = 212.55556
$ws.Range("I42").Value = 201.625
$ws.Range("K42").Value = 604.875
$ws.Range("M42").Value = -374.875
$ws.Range("H76").Value = 7995
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = $null
$ws.Range("H79").Value = 7995
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = $null
$ws.Range("H86").Value = 11129.1
$ws.Range("I86").Value = 1757.4
$ws.Range("K86").Value = 1757.4
$ws.Range("M86").Value = -634.4000000000001
$ws.Range("H89").Value = 11129.1
$ws.Range("I89").Value = 1757.4
$ws.Range("K89").Value = 8787
$ws.Range("M89").Value = -3171
$ws.Range("H100").Value = 9999
$ws.Range("I100").Value = 9999
$ws.Range("K100").Value = 9999
$ws.Range("M100").Value = -9458

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 7196.3335
$ws.Range("J15").Value = 7196.3335
$ws.Range("L15").Value = 7196.3335
$ws.Range("N15").Value = -7896.3335
$ws.Range("H16").Value = 22815.6
$ws.Range("I16").Value = 25519.5
$ws.Range("J16").Value = 12000
$ws.Range("K16").Value = 25519.5
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = -25232.5
$ws.Range("N16").Value = -12574
$ws.Range("H110").Value = 4311.769
$ws.Range("I110").Value = 2763.5
$ws.Range("K110").Value = 2763.5
$ws.Range("M110").Value = -718.5
$ws.Range("H122").Value = 2836.6428
$ws.Range("I122").Value = 2836.6428
$ws.Range("K122").Value = 8509.928400000001
$ws.Range("M122").Value = -6059.928400000001
$ws.Range("H132").Value = 3337.7585
$ws.Range("I132").Value = 3064.1072
$ws.Range("J132").Value = 11000
$ws.Range("K132").Value = 9192.321599999999
$ws.Range("L132").Value = 33000
$ws.Range("M132").Value = -6662.321599999999
$ws.Range("N132").Value = -38060

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1022
$ws.Range("I64").Value = 1228.2
$ws.Range("J64").Value = 792.8889
$ws.Range("K64").Value = 1228.2
$ws.Range("L64").Value = 792.8889
$ws.Range("M64").Value = -1003.2
$ws.Range("N64").Value = -1242.8889
$ws.Range("H67").Value = 1022
$ws.Range("I67").Value = 1228.2
$ws.Range("J67").Value = 792.8889
$ws.Range("K67").Value = 1228.2
$ws.Range("L67").Value = 792.8889
$ws.Range("M67").Value = -448.2
$ws.Range("N67").Value = -2352.8889
$ws.Range("H107").Value = 1348.125
$ws.Range("I107").Value = 705.8333
$ws.Range("J107").Value = 3275
$ws.Range("K107").Value = 705.8333
$ws.Range("L107").Value = 3275
$ws.Range("M107").Value = 1214.1667
$ws.Range("N107").Value = -7115

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 106.52941
$ws.Range("I7").Value = 120
$ws.Range("J7").Value = 74.2
$ws.Range("K7").Value = 120
$ws.Range("L7").Value = 74.2
$ws.Range("M7").Value = -7
$ws.Range("N7").Value = -300.2
$ws.Range("H120").Value = 21999.4
$ws.Range("J120").Value = 25000
$ws.Range("L120").Value = 25000
$ws.Range("N120").Value = -32258
$ws.Range("H132").Value = 3257.8
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null
$ws.Range("H134").Value = 1957.5
$ws.Range("I134").Value = 2165.375
$ws.Range("J134").Value = 1126
$ws.Range("K134").Value = 6496.125
$ws.Range("L134").Value = 3378
$ws.Range("M134").Value = -3961.125
$ws.Range("N134").Value = -8448

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 199.07692
$ws.Range("I12").Value = 47
$ws.Range("J12").Value = 244.7
$ws.Range("K12").Value = 141
$ws.Range("L12").Value = 734.0999999999999
$ws.Range("M12").Value = 32
$ws.Range("N12").Value = -1080.1
$ws.Range("H34").Value = 3055
$ws.Range("I34").Value = 466.33334
$ws.Range("J34").Value = 4025.75
$ws.Range("K34").Value = 1399.00002
$ws.Range("L34").Value = 12077.25
$ws.Range("M34").Value = -1315.00002
$ws.Range("N34").Value = -12245.25
$ws.Range("H38").Value = 37.090908
$ws.Range("I38").Value = 38.857143
$ws.Range("K38").Value = 116.571429
$ws.Range("M38").Value = 230.428571
$ws.Range("H39").Value = 3376.6
$ws.Range("J39").Value = 3628.3333
$ws.Range("L39").Value = 10884.9999
$ws.Range("N39").Value = -11472.9999
$ws.Range("H55").Value = 5786.1177
$ws.Range("I55").Value = 699.8333
$ws.Range("J55").Value = 8560.454
$ws.Range("K55").Value = 2099.4999
$ws.Range("L55").Value = 25681.362
$ws.Range("M55").Value = -1922.4999
$ws.Range("N55").Value = -26035.362
$ws.Range("H109").Value = 3380
$ws.Range("I109").Value = 915
$ws.Range("J109").Value = 6666.6665
$ws.Range("K109").Value = 2745
$ws.Range("L109").Value = 19999.9995
$ws.Range("M109").Value = -1705
$ws.Range("N109").Value = -22079.9995
$ws.Range("H115").Value = 700
$ws.Range("I115").Value = 1000
$ws.Range("J115").Value = 550
$ws.Range("K115").Value = 3000
$ws.Range("L115").Value = 1650
$ws.Range("M115").Value = -1825
$ws.Range("N115").Value = -4000
$ws.Range("H138").Value = 1447.9
$ws.Range("I138").Value = 620.5
$ws.Range("J138").Value = 1999.5
$ws.Range("K138").Value = 1861.5
$ws.Range("L138").Value = 5998.5
$ws.Range("M138").Value = 3278.5
$ws.Range("N138").Value = -16278.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2453.7144
$ws.Range("I102").Value = 2453.7144
$ws.Range("K102").Value = 2453.7144
$ws.Range("M102").Value = -831.7143999999998
$ws.Range("H122").Value = 43923.117
$ws.Range("I122").Value = 43335.285
$ws.Range("K122").Value = 130005.855
$ws.Range("M122").Value = -127555.855
$ws.Range("H141").Value = 35000
$ws.Range("J141").Value = 35000
$ws.Range("L141").Value = 35000
$ws.Range("N141").Value = -45360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1169.4762
$ws.Range("I22").Value = 1005.0769
$ws.Range("K22").Value = 1005.0769
$ws.Range("M22").Value = -710.0769
$ws.Range("H27").Value = 1169.4762
$ws.Range("I27").Value = 1005.0769
$ws.Range("K27").Value = 1005.0769
$ws.Range("M27").Value = -898.0769
$ws.Range("H55").Value = 746.1429000000001
$ws.Range("I55").Value = 722.36365
$ws.Range("K55").Value = 722.36365
$ws.Range("M55").Value = -549.36365
$ws.Range("H132").Value = 15419.038
$ws.Range("J132").Value = 11441.777
$ws.Range("L132").Value = 34325.331
$ws.Range("N132").Value = -39385.331

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 14000
$ws.Range("J43").Value = 14000
$ws.Range("L43").Value = 14000
$ws.Range("N43").Value = -14298
$ws.Range("H132").Value = 437
$ws.Range("I132").Value = 437
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1311
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 1219
$ws.Range("N132").Value = $null
$ws.Range("H135").Value = 46333
$ws.Range("I135").Value = 18999.5
$ws.Range("K135").Value = 18999.5
$ws.Range("M135").Value = -13929.5
$ws.Range("H140").Value = 73000
$ws.Range("J140").Value = 73000
$ws.Range("L140").Value = 73000
$ws.Range("N140").Value = -83360
$ws.Range("H141").Value = 59999.668
$ws.Range("J141").Value = 59999.668
$ws.Range("L141").Value = 59999.668
$ws.Range("N141").Value = -70359.66800000001
